$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits at the
#    start of the third paragraph, right before the picture run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Reposition the second picture (in the third paragraph) - the
#    anchor's horizontal/vertical offsets changed. EMU -> points
#    (1 pt = 12700 EMU).
$shape = $d.Shapes(2)
$shape.Left = -190230 / 12700
$shape.Top = 274546 / 12700

# 3. Re-insert the "_GoBack" bookmark after "Screenshot of ROS ma",
#    splitting the trailing text run into "Screenshot of ROS ma" and
#    "ster listening" around it (mirrors where the cursor was left
#    after the last edit).
$para3 = $d.Paragraphs(3).Range
$splitAt = $para3.Start + "Screenshot of ROS ma".Length
$bmRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $bmRange)
